# "update for remote CNOT" - refresh the lookahead weighting result row (row 8,
# "qft_10") now that remote CNOT is accounted for.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 ("qft_10"): new remote-CNOT run used a 4-weight lookahead list and a
# slightly lower gate count.
$ws.Range("I8").Value = 536
$ws.Range("J8").Value = "[1, 0.8, 0.6, 0.4]"

# Column J needs to widen to fit the longer weighting string.
$ws.Columns.Item(10).ColumnWidth = 13.571428571428571

# Move the active selection to K12 to match the author's saved view.
$ws.Range("K12").Select()

# The author also minimized the window before saving; best-effort (window
# chrome is not always persisted by every host, so ignore failures here).
try {
    $excel.ActiveWindow.WindowState = -4140
} catch {
}
